# "Add files via upload" – update the Functional Requirements list.
#
# The cell that held requirement F 4.3 ("... untuk melakukan pembayaran")
# is updated to describe the booking step of the reservation transaction
# instead ("... untuk melakukan booking pada transaksi reservasi").
#
# Everything else in the sheet (headers, merged cells, other requirement
# rows) stays as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C11").Value = "F 4.3 Sistem memungkinkan pengguna yang telah memesan ruangan untuk melakukan booking pada transaksi reservasi"

# Restore the cursor/selection to where the author left it (row 12, col C),
# matching the saved view state in the workbook.
$ws.Activate() | Out-Null
$ws.Range("C12").Select() | Out-Null
